# Updated cryptos list with GitHub Actions scraper refresh.
# Rewrites the Price (D) / Volume(1h) (E) columns with the latest scraped
# values, and fixes the PaxDollar/EnergySwap row ordering (rows 47-48).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price-column cells whose new text looks like a plain number get an
# explicit Text number format first so COM stores them as strings
# (preserving exact digits/trailing zeros) instead of auto-converting
# them to floating point numbers.

$ws.Range("D2").Value = '27.091.05'
$ws.Range("E2").Value = '  +0.66%  '
$ws.Range("D3").Value = '1.890.84'
$ws.Range("E3").Value = '  +1.50%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.81'
$ws.Range("E5").Value = '  +0.62%  '
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5150'
$ws.Range("E7").Value = '  +1.80%  '
$ws.Range("E8").Value = '  +3.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07196'
$ws.Range("E9").Value = '  +0.34%  '
$ws.Range("E10").Value = '  +1.49%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9018'
$ws.Range("E11").Value = '  +0.74%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07661'
$ws.Range("E12").Value = '  +2.30%  '
$ws.Range("D13").Value = '1.889.98'
$ws.Range("E13").Value = '  +1.61%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '94.24'
$ws.Range("E14").Value = '  +1.78%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.244'
$ws.Range("E15").Value = '  +0.24%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.002'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008496'
$ws.Range("E17").Value = '  +0.08%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.44'
$ws.Range("E18").Value = '  +1.85%  '
$ws.Range("E19").Value = '  +0.14%  '
$ws.Range("D20").Value = '27.127.39'
$ws.Range("E20").Value = '  +0.66%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.060'
$ws.Range("E21").Value = '  +0.45%  '
$ws.Range("D22").Value = '2.138.79'
$ws.Range("E22").Value = '  +3.56%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.57'
$ws.Range("E23").Value = '  +1.86%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.384'
$ws.Range("E24").Value = '  -0.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.290'
$ws.Range("E25").Value = '  +10.39%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '146.62'
$ws.Range("E26").Value = '  -0.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.746'
$ws.Range("E27").Value = '  -2.43%  '
$ws.Range("E28").Value = '  +0.98%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '114.48'
$ws.Range("E29").Value = '  +1.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.922'
$ws.Range("E30").Value = '  +5.14%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.801'
$ws.Range("E31").Value = '  +2.04%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09201'
$ws.Range("E32").Value = '  -0.53%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05049'
$ws.Range("E33").Value = '  -1.07%  '
$ws.Range("E34").Value = '  +7.26%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7689'
$ws.Range("E35").Value = '  +1.84%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.996'
$ws.Range("E36").Value = '  +0.61%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.285'
$ws.Range("E37").Value = '  +0.46%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.593'
$ws.Range("E38").Value = '  +2.09%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5601'
$ws.Range("E39").Value = '  +1.08%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01990'
$ws.Range("E40").Value = '  -0.35%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.073'
$ws.Range("E41").Value = '  +0.07%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.056'
$ws.Range("E42").Value = '  +6.29%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.645'
$ws.Range("E43").Value = '  +2.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '118.38'
$ws.Range("E44").Value = '  -0.28%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1501'
$ws.Range("E45").Value = '  +2.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4811'
$ws.Range("E46").Value = '  +2.53%  '
$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.001'
$ws.Range("E47").Value = '  +0.17%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.002'
$ws.Range("E48").Value = '  +0.41%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.599'
$ws.Range("E49").Value = '  +2.36%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '37.59'
$ws.Range("E50").Value = '  +1.98%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '64.01'
$ws.Range("E51").Value = '  +1.69%  '
